$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 1/2: add the new "Quantity" header + first data row's quantity.
# (shared string index 5 = "Quantity")
# ---------------------------------------------------------------------------
$ws.Range("D1").Value = "Quantity"
$ws.Range("D2").Value = 1

# ---------------------------------------------------------------------------
# New rows 3-10. Cell values are written in a specific order so that the
# shared-string table is built up in the same sequence as the original
# authoring session (matches the target workbook's sharedStrings.xml order).
# ---------------------------------------------------------------------------

# Row 3
$ws.Range("A3").Value = "08055C104JAT2A"
# Row 4
$ws.Range("A4").Value = "08055A470JAT2A"
$ws.Range("C4").Value = "47pF"
# Row 5
$ws.Range("A5").Value = "ERJ-P06J472V"
$ws.Range("C5").Value = "4.7k" + [char]0x03A9
$ws.Range("C5").Characters(5, 1).Font.Name = "Calibri"
$ws.Range("C5").Characters(5, 1).Font.Size = 11
# Row 6
$ws.Range("A6").Value = "ERJ-P06J103V"
$ws.Range("C6").Value = "10k" + [char]0x03A9
$ws.Range("C6").Characters(4, 1).Font.Name = "Calibri"
$ws.Range("C6").Characters(4, 1).Font.Size = 11
# Row 3 description (added after rows 4-6 were drafted, matching author order)
$ws.Range("C3").Value = "0.1uF"
# Row 7
$ws.Range("A7").Value = "10118194-0001LF"
$ws.Range("C7").Value = "USB mini B"
# Row 8
$ws.Range("A8").Value = "ERJ-6GEYJ271V"
$ws.Range("C8").Value = "270" + [char]0x03A9
$ws.Range("C8").Characters(4, 1).Font.Name = "Calibri"
$ws.Range("C8").Characters(4, 1).Font.Size = 11
# Row 9
$ws.Range("A9").Value = "5988170107F"
$ws.Range("C9").Value = "Green LED"
# Row 10
$ws.Range("A10").Value = "ERJ-6GEY0R00V"
$ws.Range("C10").Value = "0" + [char]0x03A9
$ws.Range("C10").Characters(2, 1).Font.Name = "Calibri"
$ws.Range("C10").Characters(2, 1).Font.Size = 11

# ---------------------------------------------------------------------------
# Quantity (column D) values for the new rows.
# ---------------------------------------------------------------------------
$ws.Range("D3").Value = 2
$ws.Range("D4").Value = 2
$ws.Range("D5").Value = 1
$ws.Range("D6").Value = 1
$ws.Range("D7").Value = 1
$ws.Range("D8").Value = 2
$ws.Range("D9").Value = 2
$ws.Range("D10").Value = 1

# ---------------------------------------------------------------------------
# Column B: Digikey HYPERLINK formulas (same style as the existing B2 cell).
# ---------------------------------------------------------------------------
$ws.Range("B3").Formula = '=HYPERLINK("http://www.digikey.ca/product-detail/en/avx-corporation/08055C104JAT2A/478-3352-1-ND/930144","Digikey - 478-3352-1-ND")'
$ws.Range("B3").Style = "Hyperlink"

$ws.Range("B4").Formula = '=HYPERLINK("http://www.digikey.ca/product-detail/en/avx-corporation/08055A470JAT2A/478-1312-1-ND/564344","Digikey - 478-1312-1-ND")'
$ws.Range("B4").Style = "Hyperlink"

$ws.Range("B5").Formula = '=HYPERLINK("http://www.digikey.ca/product-detail/en/panasonic-electronic-components/ERJ-P06J472V/P4.7KADCT-ND/525517","Digikey -  P4.7KADCT-ND")'
$ws.Range("B5").Style = "Hyperlink"

$ws.Range("B6").Formula = '=HYPERLINK("http://www.digikey.ca/product-detail/en/panasonic-electronic-components/ERJ-P06J103V/P10KADCT-ND/525438","Digikey - P10KADCT-ND")'
$ws.Range("B6").Style = "Hyperlink"

$ws.Range("B7").Formula = '=HYPERLINK("http://www.digikey.ca/product-detail/en/amphenol-fci/10118194-0001LF/609-4618-1-ND/2785382","Digikey - 609-4618-1-ND")'
$ws.Range("B7").Style = "Hyperlink"

$ws.Range("B8").Formula = '=HYPERLINK("http://www.digikey.ca/product-search/en?keywords=P270ACT-ND","Digikey - P270ACT-ND")'
$ws.Range("B8").Style = "Hyperlink"

$ws.Range("B9").Formula = '=HYPERLINK("http://www.digikey.ca/product-detail/en/dialight/5988170107F/350-2044-1-ND/1291348","Digikey - 350-2044-1-ND")'
$ws.Range("B9").Style = "Hyperlink"

$ws.Range("B10").Formula = '=HYPERLINK("http://www.digikey.ca/product-detail/en/panasonic-electronic-components/ERJ-6GEY0R00V/P0.0ACT-ND/82955","Digikey - P0.0ACT-ND")'
$ws.Range("B10").Style = "Hyperlink"

# ---------------------------------------------------------------------------
# Column A style for the new part-number cells: small Arial font, black text,
# vertically centred and word-wrapped. Built once on A3, then propagated to
# A4:A10 via a formats-only paste so every cell shares a single style record.
# ---------------------------------------------------------------------------
$ws.Range("A3").Font.Color = 0
$ws.Range("A3").Font.Size = 7
$ws.Range("A3").Font.Name = "Arial"
$ws.Range("A3").VerticalAlignment = -4108
$ws.Range("A3").WrapText = $true

$ws.Range("A3").Copy()
$ws.Range("A4:A10").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# ---------------------------------------------------------------------------
# Sheet view / print setup to match the final authored state.
# ---------------------------------------------------------------------------
$ws.PageSetup.Orientation = 1
$ws.Range("B18").Select()
